$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 6 ("Content Placeholder 2", shape 20)
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(20)
$tr6 = $shp6.TextFrame.TextRange

# Paragraph 4 is two runs:
#   "MSR6: spec minimum use-cases, architecture, YANG"
#   " spec, ? early inception of drafts for the other WG (pass over) ?"
# Only the second run's text changes, so target it precisely via Characters()
# (start/length of just that run) to avoid touching / re-splitting the first
# (unchanged) run.
$para4 = $tr6.Paragraphs(4)
$prefixLen = "MSR6: spec minimum use-cases, architecture, YANG".Length
$run2 = $tr6.Characters($para4.Start + $prefixLen, $para4.Length - $prefixLen)
$run2.Text = " spec, ? Pass over to responsible W when minimum quality met ?!"

# Paragraph 7 is two runs:
#   "Additional reality check with 6MAN"
#   " "
# Only the first run's text changes; target only that run's character span.
$para7 = $tr6.Paragraphs(7)
$run1Len = "Additional reality check with 6MAN".Length
$run1 = $tr6.Characters($para7.Start, $run1Len)
$run1.Text = "Additional reality check with V6OPS"

# ---------------------------------------------------------------------------
# Slide 7 ("Title 1", shape 1 and "Content Placeholder 2", shape 2)
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)

# Title: single run. Replace the whole paragraph span (via Characters, not
# Paragraphs().Text) so the engine doesn't fragment it into multiple runs
# when old/new text share a common prefix/suffix.
$title7 = $s7.Shapes.Item(1)
$trTitle = $title7.TextFrame.TextRange
$paraTitle = $trTitle.Paragraphs(1)
$spanTitle = $trTitle.Characters($paraTitle.Start, $paraTitle.Length)
$spanTitle.Text = "SP -> DCN: Build once, sell twice ?!!"

$shp7 = $s7.Shapes.Item(2)
$tr7 = $shp7.TextFrame.TextRange

# Paragraph 3: single run.
$para3 = $tr7.Paragraphs(3)
$span3 = $tr7.Characters($para3.Start, $para3.Length)
$span3.Text = "SRv6/SRH less necessary for TE (FlowLabel because of ECMP etc..)."

# Paragraph 4 was originally split over two runs:
#   "SRv6 N"
#   "ot needed for MSR6 if we specify appropriately !"
# Target text collapses this into a single run, so replace the whole
# paragraph span (via Characters, not Paragraphs().Text) which merges the
# two runs into one while keeping the (empty) run-properties intact.
$para4b = $tr7.Paragraphs(4)
$span4 = $tr7.Characters($para4b.Start, $para4b.Length)
$span4.Text = "SRv6 may just rely on Destination Address SID semantics without SRH."

# Paragraph 7: single run. Old/new text share a prefix+suffix
# ("Additional work for MSR6 " / "?"), which would otherwise get split into
# 3 runs by a naive Paragraphs().Text assignment; use the whole-span
# Characters() replace instead to keep it a single run.
$para7b = $tr7.Paragraphs(7)
$span7 = $tr7.Characters($para7b.Start, $para7b.Length)
$span7.Text = "Additional work for MSR6 in DCN ?"
